$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 634. This shifts the former rows
# 634:697 down to 635:698 (carrying their values/formatting along),
# growing the used range from A1:R697 to A1:R698.
$ws.Rows.Item(634).Insert()

# Populate the freshly inserted (blank) row 634 with the new weekly
# price-report record.
$ws.Range("A634").Value = 10
$ws.Range("B634").Value = "Vega Modelo de Temuco"
$ws.Range("C634").Value = "La Araucanía"
$ws.Range("D634").Value = 45166
$ws.Range("E634").Value = 9
$ws.Range("F634").Value = 100112037
$ws.Range("G634").Value = "Cebollín"
$ws.Range("H634").Value = "Sin especificar"
$ws.Range("I634").Value = "Primera"
$ws.Range("J634").Value = 80
$ws.Range("K634").Value = 8000
$ws.Range("L634").Value = 8000
$ws.Range("M634").Value = 8000
$ws.Range("N634").Value = "$/docena de paquetes"
$ws.Range("O634").Value = "Provincia de Cautín"
$ws.Range("P634").Value = 667
$ws.Range("Q634").Value = 12
$ws.Range("R634").Value = "Hortaliza"
